# Update "想去人数" (number of interested attendees) values in column F
# for the "展览" sheet and the aggregated "全部类型" sheet, matching the
# newly generated gh-pages data output.

$wb = $excel.ActiveWorkbook

# Sheet "展览" - rows keyed by their F-column cell
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 693
$wsExhibit.Range("F5").Value = 580
$wsExhibit.Range("F7").Value = 2838
$wsExhibit.Range("F9").Value = 8092
$wsExhibit.Range("F11").Value = 478
$wsExhibit.Range("F12").Value = 51
$wsExhibit.Range("F13").Value = 419
$wsExhibit.Range("F14").Value = 50

# Sheet "全部类型" - same events, different row offsets
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 693
$wsAll.Range("F5").Value = 580
$wsAll.Range("F9").Value = 2838
$wsAll.Range("F11").Value = 8092
$wsAll.Range("F13").Value = 478
$wsAll.Range("F14").Value = 51
$wsAll.Range("F17").Value = 419
$wsAll.Range("F18").Value = 50
